# issue #5: stock data from json to db
#
# The "股票" (stock) worksheet gains three new columns:
#   - "category" is inserted right after "property_category" with value
#     "normal" on every data row (it mirrors the output/normal export
#     folder the workbook file lives in),
#   - "source_file" is appended at the end with value "tmp93201" (the name
#     of the temp source file the row was generated from),
#   - "index" is appended at the end, copying each row's original
#     identifier that is already stored in column A.
#
# Concretely, the existing columns I (date), J (legislator_name) and
# K (legislator_id) each shift one column to the right (I->J, J->K, K->L),
# the brand-new "category" column becomes I, and two more brand new
# columns, "source_file" (M) and "index" (N), are appended after L.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$lastRow = 23
$headerRow = 1

# Work right-to-left, column by column, so a later Copy never clobbers data
# we still need to read. Each Copy() also carries over the source cell's
# formatting (style), which we then refresh with the intended value.

# 1) K (legislator_id) -> L
for ($r = $headerRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 11).Copy($ws.Cells.Item($r, 12))
}

# 2) J (legislator_name) -> K
for ($r = $headerRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 10).Copy($ws.Cells.Item($r, 11))
}

# 3) I (date) -> J
for ($r = $headerRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Copy($ws.Cells.Item($r, 10))
}

# 4) New "category" column header + values in I (reuse H's style, the
#    neighbouring "property_category" column, for consistent formatting)
for ($r = $headerRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Copy($ws.Cells.Item($r, 9))
}
$ws.Cells.Item($headerRow, 9).Value = "category"
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = "normal"
}

# 5) New "source_file" column (M) header + values (reuse L's style)
for ($r = $headerRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 12).Copy($ws.Cells.Item($r, 13))
}
$ws.Cells.Item($headerRow, 13).Value = "source_file"
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 13).Value = "tmp93201"
}

# 6) New "index" column (N) header + values (copy of column A, the row id;
#    reuse H's style so data rows keep the plain "s=2" look, not A's bold
#    "s=1" id style)
for ($r = $headerRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Copy($ws.Cells.Item($r, 14))
}
$ws.Cells.Item($headerRow, 14).Value = "index"
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 14).Value = $ws.Cells.Item($r, 1).Value()
}
